$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Localización" column (old column D, value "36S78W") is being split
# into two new columns: "Latitud" and "Longitud". Insert a new column
# before D so the old D ("Localización"/"36S78W") and E ("Tipo"/1) shift
# right to E/F, then overwrite D/E with the new header + numeric values.
$ws.Columns("D").Insert()

$ws.Range("D1").Value = "Latitud"
$ws.Range("E1").Value = "Longitud"

$ws.Range("D3").Value = 12.569
$ws.Range("E3").Value = 156.15

# Match the new selection recorded in the workbook.
[void]$ws.Range("D3").Select()
